# Week 10 and Women Tell All Update
# Update the "Contestants" sheet with this week's new data:
#  - Women Tell All (WTL, column K) points for several contestants
#  - Rachel's Ep10 (column L) points and her new "Eliminated" status (column M)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contestants")

# Women Tell All (WTL) points
$ws.Range("K3").Value  = 5    # Christen
$ws.Range("K4").Value  = 45   # Corinne
$ws.Range("K5").Value  = 45   # Danielle L.
$ws.Range("K6").Value  = 5    # Danielle M.
$ws.Range("K8").Value  = 5    # Hailey
$ws.Range("K10").Value = 45   # Kristina
$ws.Range("K12").Value = 30   # Rachel
$ws.Range("K13").Value = 5    # Sarah

# Rachel's week 10 points and updated status
$ws.Range("L12").Value = 130
$ws.Range("M12").Value = "Eliminated"

# Reflect the cells the author was last working in/around
$ws.Range("D1:J13").EntireColumn.ColumnWidth = 9.14
[void]$ws.Range("B4").Select()
